# Actualización automática del tracker
# Rellena las columnas "resultado" (G) y "profit" (H) de las filas
# que aún estaban vacías, según el resultado real de cada evento.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Resultado = "Fallo";   Profit = -1 },
    @{ Row = 3;  Resultado = "Acierto"; Profit = 0.83 },
    @{ Row = 4;  Resultado = "Fallo";   Profit = -1 },
    @{ Row = 8;  Resultado = "Fallo";   Profit = -1 },
    @{ Row = 12; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 17; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 21; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 29; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 35; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 37; Resultado = "Fallo";   Profit = -1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}
